# Dev IV Project Rubric - mark milestone-I features as complete ("X")
# for the rows that were actually achieved; per commit message
# "Skybox known non-working, must fix later" the Infinite Sky Box row
# is reverted back to not-yet-graded.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Rows whose "Milestone Completed(X)" column (F) should be marked X.
# Every one of these rows is tagged Milestone "I" in column E.
$xRows = @(4,5,6,18,21,22,30,31,32,33,34,35,36,38,39,40,55,65,66)

# Make sure column E carries the "I" milestone tag for every one of those
# rows *before* column F is touched (rows 34/36 didn't have it yet) - the
# G-column EXACT()/IF() formulas must see both inputs already in place so
# they recompute correctly on the same pass.
foreach ($r in $xRows) {
    $ws.Range("E$r").Value = "I"
}
foreach ($r in $xRows) {
    $ws.Range("F$r").Value = "X"
}

# Row 23 (Infinite Sky Box) is known non-working - clear its milestone tag.
$ws.Range("E23").Value = ""

# Carry-over rows at the bottom of the sheet (Effective use of GIT / API cleanup)
$ws.Range("C90").Value = "X"
$ws.Range("C91").Value = "X"

# Make sure every dependent formula (milestone totals, carry-over, etc.)
# picks up the new marks before the workbook is saved.
$excel.CalculateFull()

# Restore the view/selection state recorded for this save.
$ws.Range("F34").Select()
